$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 939.4
$ws.Cells.Item(4, 9).Value = 861.75
$ws.Cells.Item(4, 11).Value = 861.75
$ws.Cells.Item(4, 13).Value = -747.75
$ws.Cells.Item(19, 8).Value = 1402.7142
$ws.Cells.Item(19, 9).Value = 1599.5
$ws.Cells.Item(19, 11).Value = 1599.5
$ws.Cells.Item(19, 13).Value = -1424.5
$ws.Cells.Item(28, 8).Value = 821.8889
$ws.Cells.Item(28, 9).Value = 798.875
$ws.Cells.Item(28, 11).Value = 798.875
$ws.Cells.Item(28, 13).Value = -313.875
$ws.Cells.Item(31, 8).Value = 669.1429000000001
$ws.Cells.Item(31, 9).Value = 669.1429000000001
$ws.Cells.Item(31, 11).Value = 2007.4287
$ws.Cells.Item(31, 13).Value = -1777.4287
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()
$ws.Cells.Item(94, 8).Value = 662.5
$ws.Cells.Item(94, 9).Value = 662.5
$ws.Cells.Item(94, 11).Value = 662.5
$ws.Cells.Item(94, 13).Value = -211.5
$ws.Cells.Item(134, 8).Value = 99995
$ws.Cells.Item(134, 10).Value = 99995
$ws.Cells.Item(134, 12).Value = 99995
$ws.Cells.Item(134, 14).Value = -110135
$ws.Cells.Item(138, 8).Value = 3403.6667
$ws.Cells.Item(138, 9).Value = 2651.2856
$ws.Cells.Item(138, 10).Value = 3779.8572
$ws.Cells.Item(138, 11).Value = 7953.8568
$ws.Cells.Item(138, 12).Value = 11339.5716
$ws.Cells.Item(138, 13).Value = -2813.8568
$ws.Cells.Item(138, 14).Value = -21619.5716

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 1424.75
$ws.Cells.Item(3, 9).Value = 200
$ws.Cells.Item(3, 10).Value = 1833
$ws.Cells.Item(3, 11).Value = 200
$ws.Cells.Item(3, 12).Value = 1833
$ws.Cells.Item(3, 13).Value = -85
$ws.Cells.Item(3, 14).Value = -2063
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 13).ClearContents()
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 13).ClearContents()
$ws.Cells.Item(122, 8).Value = 2799.9
$ws.Cells.Item(122, 9).Value = 3999.5
$ws.Cells.Item(122, 10).Value = 2500
$ws.Cells.Item(122, 11).Value = 11998.5
$ws.Cells.Item(122, 12).Value = 7500
$ws.Cells.Item(122, 13).Value = -9548.5
$ws.Cells.Item(122, 14).Value = -12400

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 14).ClearContents()
$ws.Cells.Item(31, 8).Value = 1500
$ws.Cells.Item(31, 10).Value = 1500
$ws.Cells.Item(31, 12).Value = 1500
$ws.Cells.Item(31, 14).Value = -2004
$ws.Cells.Item(86, 8).Value = 925.1429000000001
$ws.Cells.Item(86, 10).Value = 875.2
$ws.Cells.Item(86, 12).Value = 875.2
$ws.Cells.Item(86, 14).Value = -3121.2
$ws.Cells.Item(89, 8).Value = 925.1429000000001
$ws.Cells.Item(89, 10).Value = 875.2
$ws.Cells.Item(89, 12).Value = 4376
$ws.Cells.Item(89, 14).Value = -15608
$ws.Cells.Item(99, 8).Value = 7799
$ws.Cells.Item(99, 9).Value = 6665
$ws.Cells.Item(99, 10).Value = 9500
$ws.Cells.Item(99, 11).Value = 6665
$ws.Cells.Item(99, 12).Value = 9500
$ws.Cells.Item(99, 13).Value = -5167
$ws.Cells.Item(99, 14).Value = -12496
$ws.Cells.Item(107, 8).Value = 2211
$ws.Cells.Item(107, 9).Value = 2211
$ws.Cells.Item(107, 11).Value = 2211
$ws.Cells.Item(107, 13).Value = -291

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(14, 8).Value = 516.6667
$ws.Cells.Item(14, 10).Value = 516.6667
$ws.Cells.Item(14, 12).Value = 516.6667
$ws.Cells.Item(14, 14).Value = -856.6667
$ws.Cells.Item(22, 8).Value = 4849.5
$ws.Cells.Item(22, 9).Value = 4849.5
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 4849.5
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -4499.5
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 4375
$ws.Cells.Item(86, 9).Value = 4900
$ws.Cells.Item(86, 11).Value = 4900
$ws.Cells.Item(86, 13).Value = -3777
$ws.Cells.Item(89, 8).Value = 4375
$ws.Cells.Item(89, 9).Value = 4900
$ws.Cells.Item(89, 11).Value = 24500
$ws.Cells.Item(89, 13).Value = -18884
$ws.Cells.Item(99, 8).Value = 1799.5
$ws.Cells.Item(99, 9).Value = 1799.5
$ws.Cells.Item(99, 11).Value = 1799.5
$ws.Cells.Item(99, 13).Value = -301.5
$ws.Cells.Item(126, 8).Value = 1799.5
$ws.Cells.Item(126, 9).Value = 1799.5
$ws.Cells.Item(126, 11).Value = 5398.5
$ws.Cells.Item(126, 13).Value = -2928.5
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1125
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 12).Value = 3000
$ws.Cells.Item(5, 14).Value = -3224
$ws.Cells.Item(26, 8).Value = 84.75
$ws.Cells.Item(26, 9).Value = 88
$ws.Cells.Item(26, 10).Value = 75
$ws.Cells.Item(26, 11).Value = 264
$ws.Cells.Item(26, 12).Value = 225
$ws.Cells.Item(26, 13).Value = 24
$ws.Cells.Item(26, 14).Value = -801
$ws.Cells.Item(75, 8).Value = 1595
$ws.Cells.Item(75, 10).Value = 1595
$ws.Cells.Item(75, 12).Value = 4785
$ws.Cells.Item(75, 14).Value = -6781
$ws.Cells.Item(78, 8).Value = 1595
$ws.Cells.Item(78, 10).Value = 1595
$ws.Cells.Item(78, 12).Value = 14355
$ws.Cells.Item(78, 14).Value = -24339
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 14).ClearContents()
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 11).Value = 0
$ws.Cells.Item(109, 13).ClearContents()
$ws.Cells.Item(115, 8).Value = 4999
$ws.Cells.Item(115, 9).Value = 4999
$ws.Cells.Item(115, 11).Value = 14997
$ws.Cells.Item(115, 13).Value = -13822
$ws.Cells.Item(135, 8).Value = 1125
$ws.Cells.Item(135, 10).Value = 1000
$ws.Cells.Item(135, 12).Value = 9000
$ws.Cells.Item(135, 14).Value = -14070
$ws.Cells.Item(140, 8).Value = 558.8182
$ws.Cells.Item(140, 9).Value = 558.8182
$ws.Cells.Item(140, 11).Value = 1676.4546
$ws.Cells.Item(140, 13).Value = 3503.5454

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 14).ClearContents()
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).ClearContents()
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).ClearContents()
$ws.Cells.Item(97, 8).Value = 1513.1428
$ws.Cells.Item(97, 9).Value = 2149.75
$ws.Cells.Item(97, 10).Value = 664.3333
$ws.Cells.Item(97, 11).Value = 2149.75
$ws.Cells.Item(97, 12).Value = 664.3333
$ws.Cells.Item(97, 13).Value = -1653.75
$ws.Cells.Item(97, 14).Value = -1656.3333
$ws.Cells.Item(99, 8).Value = 2999
$ws.Cells.Item(99, 9).Value = 2999
$ws.Cells.Item(99, 11).Value = 2999
$ws.Cells.Item(99, 13).Value = -753
$ws.Cells.Item(107, 8).Value = 551.5
$ws.Cells.Item(113, 8).Value = 2000
$ws.Cells.Item(113, 9).Value = 2000
$ws.Cells.Item(113, 11).Value = 2000
$ws.Cells.Item(113, 13).Value = 170
$ws.Cells.Item(122, 8).Value = 1006.6667
$ws.Cells.Item(122, 9).Value = 1110
$ws.Cells.Item(122, 11).Value = 3330
$ws.Cells.Item(122, 13).Value = -880
$ws.Cells.Item(126, 8).Value = 5399.4
$ws.Cells.Item(126, 9).Value = 5399.4
$ws.Cells.Item(126, 11).Value = 16198.2
$ws.Cells.Item(126, 13).Value = -13728.2
$ws.Cells.Item(134, 8).Value = 99998.336
$ws.Cells.Item(134, 10).Value = 99998.336
$ws.Cells.Item(134, 12).Value = 299995.008
$ws.Cells.Item(134, 14).Value = -305065.008

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3985
$ws.Cells.Item(46, 10).Value = 3985
$ws.Cells.Item(46, 12).Value = 3985
$ws.Cells.Item(46, 14).Value = -4361
$ws.Cells.Item(61, 8).Value = 3287.7693
$ws.Cells.Item(61, 9).Value = 3224.5
$ws.Cells.Item(61, 10).Value = 3389
$ws.Cells.Item(61, 11).Value = 3224.5
$ws.Cells.Item(61, 12).Value = 3389
$ws.Cells.Item(61, 13).Value = -3022.5
$ws.Cells.Item(61, 14).Value = -3793
$ws.Cells.Item(82, 8).Value = 1454.5834
$ws.Cells.Item(82, 9).Value = 1461.25
$ws.Cells.Item(82, 11).Value = 1461.25
$ws.Cells.Item(82, 13).Value = -1100.25
$ws.Cells.Item(85, 8).Value = 1454.5834
$ws.Cells.Item(85, 9).Value = 1461.25
$ws.Cells.Item(85, 11).Value = 1461.25
$ws.Cells.Item(85, 13).Value = -213.25
$ws.Cells.Item(100, 8).Value = 10188
$ws.Cells.Item(100, 9).Value = 3538.4
$ws.Cells.Item(100, 11).Value = 3538.4
$ws.Cells.Item(100, 13).Value = -2997.4
$ws.Cells.Item(113, 8).Value = 3287.7693
$ws.Cells.Item(113, 9).Value = 3224.5
$ws.Cells.Item(113, 10).Value = 3389
$ws.Cells.Item(113, 11).Value = 3224.5
$ws.Cells.Item(113, 12).Value = 3389
$ws.Cells.Item(113, 13).Value = -1054.5
$ws.Cells.Item(113, 14).Value = -7729
$ws.Cells.Item(132, 8).Value = 3416.3333
$ws.Cells.Item(132, 9).Value = 3879.6
$ws.Cells.Item(132, 10).Value = 1100
$ws.Cells.Item(132, 11).Value = 11638.8
$ws.Cells.Item(132, 12).Value = 3300
$ws.Cells.Item(132, 13).Value = -9108.799999999999
$ws.Cells.Item(132, 14).Value = -8360
$ws.Cells.Item(134, 8).Value = 100000
$ws.Cells.Item(134, 10).Value = 100000
$ws.Cells.Item(134, 12).Value = 100000
$ws.Cells.Item(134, 14).Value = -110140
$ws.Cells.Item(136, 8).Value = 2499
$ws.Cells.Item(136, 9).Value = 2499
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 7497
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -4947
$ws.Cells.Item(136, 14).ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(133, 8).Value = 117498.5
$ws.Cells.Item(133, 10).Value = 117498.5
$ws.Cells.Item(133, 12).Value = 117498.5
$ws.Cells.Item(133, 14).Value = -127618.5
$ws.Cells.Item(136, 8).Value = 1189.6923
$ws.Cells.Item(136, 9).Value = 1189.6923
$ws.Cells.Item(136, 11).Value = 3569.0769
$ws.Cells.Item(136, 13).Value = -1019.0769
